$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.960.16"
$ws.Range("E2").Value = "  +2.02%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.815.36"
$ws.Range("E3").Value = "  +2.43%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.49%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.32"
$ws.Range("E5").Value = "  +1.91%  "

# Row 6 - USDC
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.36%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.10%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3670"
$ws.Range("E8").Value = "  +0.09%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07258"
$ws.Range("E9").Value = "  +0.48%  "

# Row 10 - WrappedEther
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.160.64"
$ws.Range("E10").Value = "  +22.17%  "

# Row 11 - Polygon
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8641"
$ws.Range("E11").Value = "  +1.64%  "

# Row 12 - Solana
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.29"
$ws.Range("E12").Value = "  +4.86%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.409"
$ws.Range("E13").Value = "  +3.13%  "

# Row 14 - Chainlink
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.604"
$ws.Range("E14").Value = "  +2.43%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +1.56%  "

# Row 16 - Litecoin
$ws.Range("E16").Value = "  +1.88%  "

# Row 17 - BinanceUSD
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.011"
$ws.Range("E17").Value = "  +0.50%  "

# Row 18 - ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008865"
$ws.Range("E18").Value = "  +2.02%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.32%  "

# Row 20 - Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.25"
$ws.Range("E20").Value = "  +1.24%  "

# Row 21 - WrappedBTC
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.001.96"
$ws.Range("E21").Value = "  +2.16%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.186"
$ws.Range("E22").Value = "  +1.22%  "

# Row 23 - was Cosmos, now WrappedliquidstakedEther2.0 (rows 23/24 swap with updated values)
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.389.93"
$ws.Range("E23").Value = "  +19.94%  "

# Row 24 - was WrappedliquidstakedEther2.0, now Cosmos
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.00"
$ws.Range("E24").Value = "  -2.61%  "

# Row 25 - Monero
$ws.Range("E25").Value = "  +1.10%  "

# Row 26 - Toncoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.884"
$ws.Range("E26").Value = "  +1.70%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.30"
$ws.Range("E27").Value = "  +0.77%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.220"
$ws.Range("E28").Value = "  +2.67%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.898"
$ws.Range("E29").Value = "  +10.25%  "

# Row 30 - BitcoinCash
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.63"
$ws.Range("E30").Value = "  -0.10%  "

# Row 31 - Stellar
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08941"
$ws.Range("E31").Value = "  -0.16%  "

# Row 32 - ARBITRUM
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.185"
$ws.Range("E32").Value = "  +6.20%  "

# Row 33 - ImmutableX
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7461"
$ws.Range("E33").Value = "  +2.95%  "

# Row 34 - Filecoin
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.415"
$ws.Range("E34").Value = "  +1.94%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +2.31%  "

# Row 36 - Frax
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.006"
$ws.Range("E36").Value = "  +0.35%  "

# Row 37 - TrustWalletToken
$ws.Range("E37").Value = "  +4.97%  "

# Row 38 - Hedera
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05213"
$ws.Range("E38").Value = "  +0.90%  "

# Row 39 - VeChain
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01924"
$ws.Range("E39").Value = "  +1.45%  "

# Row 40 - TheSandbox
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5097"
$ws.Range("E40").Value = "  +3.30%  "

# Row 41 - Algorand
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1655"
$ws.Range("E41").Value = "  +3.03%  "

# Row 42 - MXToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.732"
$ws.Range("E42").Value = "  +8.00%  "

# Row 43 - FraxShare
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.480"
$ws.Range("E43").Value = "  +4.08%  "

# Row 44 - Aptos
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.330"
$ws.Range("E44").Value = "  +3.60%  "

# Row 45 - was PaxosStandard, now Quant (list shifted up, PaxosStandard dropped)
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.65"
$ws.Range("E45").Value = "  +1.71%  "

# Row 46 - was Quant, now EnergySwap
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.35"
$ws.Range("E46").Value = "  +2.34%  "

# Row 47 - was EnergySwap, now PaxDollar
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.006"
$ws.Range("E47").Value = "  +0.42%  "

# Row 48 - was PaxDollar, now Decentraland
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4585"
$ws.Range("E48").Value = "  +1.82%  "

# Row 49 - was Decentraland, now NEARProtocol
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.643"
$ws.Range("E49").Value = "  +3.75%  "

# Row 50 - was NEARProtocol, now Cronos
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06215"
$ws.Range("E50").Value = "  +0.27%  "

# Row 51 - was Cronos, now RenderToken (new entry)
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.834"
$ws.Range("E51").Value = "  +5.26%  "
